# Append the 2021 data row (row 11) to Sheet1, matching the yearly series that
# already runs 2012..2020 in rows 2..10. Columns E ("其他采矿业") and V ("烟草制品业")
# have no reported figure for 2021, so they are intentionally left blank -- same
# as how the other "-" placeholders already appear in the existing rows above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from the last data row (A10) down to the new row label (A11),
# then overwrite the value -- mirrors the "2021" row carrying the same style index as "2020".
$ws.Range("A10").Copy($ws.Range("A11"))

$ws.Range("A11").Value = "2021年"
$ws.Range("B11").Value = 7480.79
$ws.Range("C11").Value = 1823.18
$ws.Range("D11").Value = 287.92
$ws.Range("F11").Value = 3964.72
$ws.Range("G11").Value = 9942.33
$ws.Range("H11").Value = 868.04
$ws.Range("I11").Value = 6968.96
$ws.Range("J11").Value = 881.67
$ws.Range("K11").Value = 168071.24
$ws.Range("L11").Value = 961.09
$ws.Range("M11").Value = 172.12
$ws.Range("N11").Value = 44.57
$ws.Range("O11").Value = 1970.36
$ws.Range("P11").Value = 3075.48
$ws.Range("Q11").Value = 90.59
$ws.Range("R11").Value = 297.5
$ws.Range("S11").Value = 4173.4
$ws.Range("T11").Value = 730.1799999999999
$ws.Range("U11").Value = 23519.23
$ws.Range("W11").Value = 1297.07
$ws.Range("X11").Value = 1814.83
$ws.Range("Y11").Value = 3160.47
$ws.Range("Z11").Value = 11910.82
$ws.Range("AA11").Value = 1401.31
$ws.Range("AB11").Value = 1608.71
$ws.Range("AC11").Value = 526.1
$ws.Range("AD11").Value = 2470.73
$ws.Range("AE11").Value = 2161.1
$ws.Range("AF11").Value = 40184.8
$ws.Range("AG11").Value = 10184.62
$ws.Range("AH11").Value = 3408.22
$ws.Range("AI11").Value = 1807.94
$ws.Range("AJ11").Value = 332.95
$ws.Range("AK11").Value = 5020.53
$ws.Range("AL11").Value = 2170.26
$ws.Range("AM11").Value = 4236.83
$ws.Range("AN11").Value = 64.58
$ws.Range("AO11").Value = 3816.33
$ws.Range("AP11").Value = 3098.11
$ws.Range("AQ11").Value = 118.53
